$d = $word.ActiveDocument

# --- 1) Insert a "<FIXED>" paragraph right before the paragraph that reads
#        "- When there is an empty /n after a valid instruction, the program
#        crashes (null pointer exception)"
$find1 = $d.Content.Find
$find1.Execute("- When there is an empty /n after a valid instruction") | Out-Null
$target1Index = $find1.Parent.Paragraphs(1).Index
$target1 = $d.Paragraphs($target1Index)
$target1.Range.InsertParagraphBefore()
$d.Paragraphs($target1Index).Range.Text = "<FIXED>"

# --- 2) Insert a "<FIXED>" paragraph right before the paragraph that reads
#        "- When there is an empty space in between valid instructions"
$find2 = $d.Content.Find
$find2.Execute("- When there is an empty space") | Out-Null
$target2Index = $find2.Parent.Paragraphs(1).Index
$target2 = $d.Paragraphs($target2Index)
$target2.Range.InsertParagraphBefore()
$d.Paragraphs($target2Index).Range.Text = "<FIXED>"

# --- 3) Append two new, completely empty paragraphs at the very end of the
#        document (after the last "DSUBU R0,R0,R0" line, before the section
#        break).
$emptyParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$endRange1 = $d.Content
$endRange1.Collapse(0)
$endRange1.InsertXML($emptyParagraphXml) | Out-Null

$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertXML($emptyParagraphXml) | Out-Null
